$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0
$ws.Range("A3").Value = "Warnungen"
$ws.Range("B3").Value = 4
$ws.Range("C3").Value = 4

$ws.Range("F6").Select()
